# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.465.89"
$ws.Range("E2").Value = "  -0.21%  "

# Row 3
$ws.Range("D3").Value = "3.673.12"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'621.15"
$ws.Range("E5").Value = "  -7.64%  "

# Row 6
$ws.Range("D6").Value = "'159.09"
$ws.Range("E6").Value = "  -1.72%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("D8").Value = "'0.496"
$ws.Range("E8").Value = "  -0.44%  "

# Row 9
$ws.Range("E9").Value = "  -1.32%  "

# Row 10
$ws.Range("E10").Value = "  +1.45%  "

# Row 11
$ws.Range("D11").Value = "'0.440"
$ws.Range("E11").Value = "  -0.84%  "

# Row 12
$ws.Range("E12").Value = "  -2.62%  "

# Row 13
$ws.Range("D13").Value = "4.293.19"
$ws.Range("E13").Value = "  -0.90%  "

# Row 14
$ws.Range("D14").Value = "'32.30"
$ws.Range("E14").Value = "  -1.69%  "

# Row 15
$ws.Range("D15").Value = "3.664.23"
$ws.Range("E15").Value = "  -0.92%  "

# Row 16
$ws.Range("D16").Value = "69.470.47"
$ws.Range("E16").Value = "  -0.26%  "

# Row 17
$ws.Range("E17").Value = "  +0.61%  "

# Row 18
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'15.91"
$ws.Range("E18").Value = "  -2.50%  "

# Row 19
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'6.50"
$ws.Range("E19").Value = "  -0.36%  "

# Row 20
$ws.Range("D20").Value = "'10.29"
$ws.Range("E20").Value = "  +4.91%  "

# Row 21
$ws.Range("D21").Value = "'468.81"
$ws.Range("E21").Value = "  -1.15%  "

# Row 22
$ws.Range("E22").Value = "  -0.46%  "

# Row 23
$ws.Range("D23").Value = "'79.60"
$ws.Range("E23").Value = "  -0.98%  "

# Row 24
$ws.Range("D24").Value = "3.821.02"
$ws.Range("E24").Value = "  -0.87%  "

# Row 25
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
$ws.Range("D26").Value = "'11.13"
$ws.Range("E26").Value = "  +1.12%  "

# Row 27
$ws.Range("D27").Value = "'0.0000122"
$ws.Range("E27").Value = "  -4.83%  "

# Row 28
$ws.Range("D28").Value = "'8.66"
$ws.Range("E28").Value = "  -5.45%  "

# Row 29
$ws.Range("D29").Value = "'2.62"
$ws.Range("E29").Value = "  -2.90%  "

# Row 30
$ws.Range("D30").Value = "'1.66"
$ws.Range("E30").Value = "  -4.30%  "

# Row 31
$ws.Range("E31").Value = "  +0.31%  "

# Row 32
$ws.Range("E32").Value = "  -2.35%  "

# Row 33
$ws.Range("D33").Value = "'26.61"
$ws.Range("E33").Value = "  -1.29%  "

# Row 34
$ws.Range("D34").Value = "'6.39"
$ws.Range("E34").Value = "  -3.11%  "

# Row 35
$ws.Range("D35").Value = "3.677.04"
$ws.Range("E35").Value = "  -0.52%  "

# Row 36
$ws.Range("E36").Value = "  -3.39%  "

# Row 37
$ws.Range("D37").Value = "'8.28"
$ws.Range("E37").Value = "  -3.11%  "

# Row 38
$ws.Range("E38").Value = "  +0.02%  "

# Row 39
$ws.Range("D39").Value = "'178.71"
$ws.Range("E39").Value = "  +2.72%  "

# Row 40
$ws.Range("E40").Value = "  -0.08%  "

# Row 41
$ws.Range("D41").Value = "'2.22"
$ws.Range("E41").Value = "  -1.72%  "

# Row 42
$ws.Range("D42").Value = "'5.80"
$ws.Range("E42").Value = "  -5.40%  "

# Row 43
$ws.Range("D43").Value = "'0.0893"
$ws.Range("E43").Value = "  -2.42%  "

# Row 44
$ws.Range("D44").Value = "'0.924"
$ws.Range("E44").Value = "  -1.91%  "

# Row 45
$ws.Range("D45").Value = "'29.29"
$ws.Range("E45").Value = "  +5.50%  "

# Row 46
$ws.Range("D46").Value = "'46.67"
$ws.Range("E46").Value = "  -0.87%  "

# Row 47
$ws.Range("D47").Value = "'2.70"
$ws.Range("E47").Value = "  -2.32%  "

# Row 48
$ws.Range("D48").Value = "'7.85"
$ws.Range("E48").Value = "  -0.46%  "

# Row 49
$ws.Range("D49").Value = "'0.000264"
$ws.Range("E49").Value = "  -5.51%  "

# Row 50
$ws.Range("E50").Value = "  -5.38%  "

# Row 51
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "'0.260"
$ws.Range("E51").Value = "  -3.00%  "
